$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 5 identical rows (A1:A5); rows 1-2 carried a
# custom row height (75) and cell style (wrap text), rows 3-5 were plain.
# Target state: only 2 rows remain (A1:A2), with plain/default formatting
# and the selection sitting on A1.
#
# Deleting the two formatted rows (1-2) shifts the three plain rows
# (3-5) up to become the new rows 1-2 (keeping their default formatting),
# then we drop the extra trailing plain row so only 2 remain.
$ws.Range("A1:A2").EntireRow.Delete() | Out-Null
$ws.Range("A3").EntireRow.Delete() | Out-Null

# Park the selection on A1 (was A2).
$ws.Range("A1").Select() | Out-Null
